$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Update the validation message for test case 3 (row 4): the TaxType id
# referenced in the error message changed from 10 to 12.
$ws.Range("D4").Value = "Unable to find TaxType with id 12"

# Update the selection stored in the sheet view (no longer a single active
# cell D2, now a range selection A1:D2).
$ws.Range("A1:D2").Select()

$ws.Activate()
